# Applies the "Updated cryptos list" data refresh described by the diff:
# refreshed Price (D) and Volume(1h) (E) values for most rows, plus a full
# content swap of rows 44/45 (Stacks <-> NEARProtocol) and a replacement of
# row 51 (FlareNetwork -> Mantle).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text reads as plain text to Excel (URLs, coin names, the
# "  +x.xx%  " volume strings, and Price values containing more than one dot
# such as "51.198.31") -- no special handling required.
$plainCells = @(
    "D2",
    "E2",
    "D3",
    "E3",
    "E4",
    "E5",
    "E6",
    "E7",
    "E8",
    "E9",
    "E10",
    "E11",
    "E12",
    "D13",
    "E13",
    "E14",
    "E15",
    "D16",
    "E16",
    "E17",
    "E18",
    "D19",
    "E19",
    "E20",
    "E21",
    "D22",
    "E22",
    "E23",
    "E24",
    "E25",
    "E26",
    "E27",
    "E28",
    "E29",
    "E30",
    "E31",
    "E32",
    "E33",
    "E34",
    "E35",
    "E36",
    "E37",
    "E38",
    "E39",
    "E40",
    "E41",
    "E42",
    "E43",
    "B44",
    "C44",
    "E44",
    "B45",
    "C45",
    "E45",
    "E46",
    "E47",
    "E48",
    "D49",
    "E49",
    "E50",
    "B51",
    "C51",
    "E51"
)
$plainValues = @(
    '51.198.31',
    '  -0.81%  ',
    '3.062.87',
    '  +0.80%  ',
    '  +0.02%  ',
    '  +1.54%  ',
    '  -1.43%  ',
    '  -2.24%  ',
    '  +0.03%  ',
    '  -0.91%  ',
    '  -0.39%  ',
    '  +0.42%  ',
    '  -1.67%  ',
    '3.542.29',
    '  +0.81%  ',
    '  -2.25%  ',
    '  -1.42%  ',
    '3.066.85',
    '  +1.06%  ',
    '  +3.90%  ',
    '  -1.92%  ',
    '51.218.41',
    '  -0.85%  ',
    '  +1.60%  ',
    '  -1.41%  ',
    '0.0₃0955',
    '  -0.87%  ',
    '  -0.44%  ',
    '  -1.20%  ',
    '  -0.95%  ',
    '  -5.34%  ',
    '  +1.60%  ',
    '  +0.05%  ',
    '  -5.94%  ',
    '  -7.16%  ',
    '  -2.56%  ',
    '  +2.76%  ',
    '  +9.14%  ',
    '  +4.78%  ',
    '  +0.21%  ',
    '  -1.32%  ',
    '  -0.11%  ',
    '  -0.52%  ',
    '  +0.03%  ',
    '  +0.42%  ',
    '  -2.99%  ',
    '  -2.07%  ',
    '  -1.54%  ',
    'NEARProtocol',
    'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near',
    '  +2.86%  ',
    'Stacks',
    'https://coinranking.com/coin/mMPrMcB7+stacks-stx',
    '  -1.63%  ',
    '  -0.32%  ',
    '  +0.20%  ',
    '  -1.99%  ',
    '2.065.88',
    '  +1.58%  ',
    '  -0.16%  ',
    'Mantle',
    'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt',
    '  +9.70%  '
)
for ($i = 0; $i -lt $plainCells.Length; $i++) {
    $ws.Range($plainCells[$i]).Value = $plainValues[$i]
}

# Price cells whose new text is a "clean" decimal number (e.g. "390.30") and
# would otherwise be auto-converted to a numeric value by Excel, losing the
# trailing zero / text representation. Force text format, assign, then drop
# back to the default "Normal" style so no stray formatting remains.
$textCells = @(
    "D5",
    "D6",
    "D7",
    "D9",
    "D10",
    "D12",
    "D14",
    "D15",
    "D18",
    "D21",
    "D23",
    "D24",
    "D25",
    "D26",
    "D27",
    "D29",
    "D30",
    "D32",
    "D33",
    "D35",
    "D36",
    "D40",
    "D44",
    "D45",
    "D46",
    "D47",
    "D48",
    "D51"
)
$textValues = @(
    '390.30',
    '101.31',
    '0.532',
    '0.583',
    '36.75',
    '0.0847',
    '18.29',
    '7.66',
    '10.57',
    '12.25',
    '69.68',
    '264.06',
    '3.14',
    '7.86',
    '26.80',
    '7.12',
    '0.161',
    '10.59',
    '0.0490',
    '2.08',
    '49.91',
    '128.23',
    '3.79',
    '2.49',
    '21.61',
    '2.46',
    '2.04',
    '0.885'
)
for ($i = 0; $i -lt $textCells.Length; $i++) {
    $addr = $textCells[$i]
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $textValues[$i]
    $ws.Range($addr).Style = "Normal"
}
